$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '63.780.94'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '2.612.56'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '595.41'
$ws.Range('E5').Value = '  -1.79%  '
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.78%  '
$ws.Range('E9').Value = '  +2.01%  '
$ws.Range('E10').Value = '  +3.41%  '
$ws.Range('D11').Value = '0.386'
$ws.Range('E11').Value = '  +3.86%  '
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '27.92'
$ws.Range('E13').Value = '  +3.17%  '
$ws.Range('D14').Value = '3.083.54'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '63.642.32'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('E16').Value = '  +6.05%  '
$ws.Range('D17').Value = '2.619.33'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '12.44'
$ws.Range('E18').Value = '  +8.67%  '
$ws.Range('D19').Value = '4.72'
$ws.Range('E19').Value = '  +5.02%  '
$ws.Range('D20').Value = '349.15'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').Value = '6.89'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = '67.72'
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('E24').Value = '  +7.58%  '
$ws.Range('E25').Value = '  +4.32%  '
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = '557.00'
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').Value = '5.29'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('D35').Value = '166.65'
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  +4.05%  '
$ws.Range('D39').Value = '1.94'
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '167.00'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').Value = '39.73'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('E43').Value = '  +6.07%  '
$ws.Range('D44').Value = '0.0589'
$ws.Range('E44').Value = '  +5.04%  '
$ws.Range('D45').Value = '22.16'
$ws.Range('E45').Value = '  +2.09%  '
$ws.Range('D46').Value = '0.634'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').Value = '  +6.38%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0252'
$ws.Range('E48').Value = '  +4.41%  '
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').Value = '19.31'
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('D51').Value = '0.0₆0238'
$ws.Range('E51').Value = '  +22.41%  '
